# Fix two spelling mistakes flagged by spell-check ("err" squiggles) in the deck.
$p = $ppt.ActivePresentation

# --- Slide 1 ("Study of posibilities" title slide): "Inés" -> "Ines" -------
# Shape 2 is the subtitle ("Sous-titre 2") that starts with the author's
# first name, spelled "Inés" with an accent. Replace just that first word.
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$tr1 = $subtitle.TextFrame.TextRange

# Select the leading "Inés" (4 characters) and retype it without the accent.
# Delete + InsertBefore (rather than overwriting the Characters() range text
# in place) so the corrected word is typed fresh instead of re-using the
# old, spell-flagged run -- exactly like fixing the word in the UI.
$word1 = $tr1.Characters(1, 4)
$word1.Delete()
$tr1b = $subtitle.TextFrame.TextRange
$tr1b.InsertBefore("Ines") | Out-Null

# --- Slide 9 ("Descision table"): "Descision" -> "Decision" ---------------
$s9 = $p.Slides.Item(9)
$title9 = $s9.Shapes.Item(1)
$tr9 = $title9.TextFrame.TextRange

# "Descision" is the first 9 characters, followed by " table".
$word2 = $tr9.Characters(1, 9)
$word2.Delete()
$tr9b = $title9.TextFrame.TextRange
$rest = $tr9b.Characters(1, $tr9b.Length)
$rest.Text = "Decision table"
